$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$oldGuid = "94e35117-715b-43b5-9d4e-54bc1dad67bb"
$newGuid = "48e3ae9f-0b13-4ea0-add8-adac22dd36bb"
$oldHash = "016a76dd21889117c16de60df1eb254461145ebb"
$newHash = "485a19ee97895183dbfe095b0660fac3fb111ed3"

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$ws1.Range("A2").Value2 = "$newGuid.md"
$ws1.Range("B2").Value2 = "e2e\$newGuid.md"

foreach ($h in $ws1.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$B$2') {
        $h.TextToDisplay = "e2e\$newGuid.md"
    }
}

$ws1.Range("G2").Value2 = "2016-08-30 23:04:55"

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$ws2.Range("A2").Value2 = "$newGuid.md"

foreach ($h in $ws2.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = "$newGuid.md"
    }
}

$ws2.Range("G2").Value2 = "$newGuid.$newHash.zh-cn.xlf"
$ws2.Range("H2").Value2 = "2016-08-30 23:04:50"

foreach ($h in $ws2.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$I$2') {
        $h.Delete()
    }
}

$ws2.Range("I2").Value2 = ""
$ws2.Range("I2").Style = "Normal"
$ws2.Range("J2").Value2 = ""
$ws2.Range("K2").Value2 = "0001-01-01 00:00:00"

$ws2.Columns.Item(9).ColumnWidth = 17.8
$ws2.Columns.Item(10).ColumnWidth = 20.8

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$ws3.Range("A2").Value2 = "$newGuid.md"

foreach ($h in $ws3.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = "$newGuid.md"
    }
}

$ws3.Range("G2").Value2 = "$newGuid.$newHash.de-de.xlf"
$ws3.Range("H2").Value2 = "2016-08-30 23:04:55"

foreach ($h in $ws3.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$I$2') {
        $h.Delete()
    }
}

$ws3.Range("I2").Value2 = ""
$ws3.Range("I2").Style = "Normal"
$ws3.Range("J2").Value2 = ""
$ws3.Range("K2").Value2 = "0001-01-01 00:00:00"

$ws3.Columns.Item(9).ColumnWidth = 17.8
$ws3.Columns.Item(10).ColumnWidth = 20.8
